$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save") should match the style used by the other
# header cells (e.g. G1 "sum"), so copy that formatting over.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values for rows 2 and 3.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
